$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 — shifts existing rows 9-20 down to 10-21
# and grows the used range from A1:I20 to A1:I21.
$ws.Rows.Item(9).Insert()

# --- Populate the newly inserted row 9 ("Short point (up to 3 mtr.)") ---
$ws.Range("A9").Value = "P. point"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 97
$ws.Range("D9").Value = "'2"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "Short point (up to 3 mtr.)"
$ws.Range("F9").Value = 256
$ws.Range("G9").Value = "'24832.00"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = "'"
$ws.Range("I9").Style = "Normal"

# --- Row 8: quantity-to-date edited ---
$ws.Range("C8").Value = 14

# --- Rows 10-17 (formerly 9-16, now shifted down by the insert):
#     refreshed "Qty executed upto date" (C) and "Upto date Amount" (G) ---
$ws.Range("C10").Value = 11
$ws.Range("G10").Value = "'5192.00"
$ws.Range("G10").Style = "Normal"

$ws.Range("C11").Value = 53
$ws.Range("G11").Value = "'35086.00"
$ws.Range("G11").Style = "Normal"

$ws.Range("C12").Value = 98
$ws.Range("G12").Value = "'0.00"
$ws.Range("G12").Style = "Normal"

$ws.Range("C13").Value = 56
$ws.Range("G13").Value = "'7616.00"
$ws.Range("G13").Style = "Normal"

$ws.Range("C14").Value = 89
$ws.Range("G14").Value = "'2047.00"
$ws.Range("G14").Style = "Normal"

$ws.Range("C15").Value = 87
$ws.Range("G15").Value = "'0.00"
$ws.Range("G15").Style = "Normal"

$ws.Range("C16").Value = 11
$ws.Range("G16").Value = "'0.00"
$ws.Range("G16").Style = "Normal"

$ws.Range("C17").Value = 19
$ws.Range("G17").Value = "'0.00"
$ws.Range("G17").Style = "Normal"

# --- Summary block (rows 19 & 21, formerly 18 & 20): new Grand Total / Net Payable ---
$ws.Range("G19").Value = "'74773.00"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value = "'74773.00"
$ws.Range("H19").Style = "Normal"

$ws.Range("G21").Value = "'74773.00"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = "'74773.00"
$ws.Range("H21").Style = "Normal"
